$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 732 (the "「開かれたターと結ばれたター」" post) was removed from the
# source data. Delete that entire worksheet row; Excel will shift all rows
# below it up by one automatically.
$ws.Rows.Item(732).Delete()
